$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.700661
$ws.Range("H2").Value = 68.10198299999999
$ws.Range("I2").Value = 0.08615268874617349
$ws.Range("J2").Value = 0.08615268874617349
$ws.Range("M2").Value = 0.5373756666666667
$ws.Range("N2").Value = 1.612127
$ws.Range("O2").Value = 0.007472820128982582
$ws.Range("P2").Value = 0.007472820128982581
$ws.Range("Q2").Value = 12.198782838649
$ws.Range("R2").Value = 109.789045547841
$ws.Range("S2").Value = 0.0006438035466283764
$ws.Range("T2").Value = 0.0006438035466283763
$ws.Range("G3").Value = 22.700661
$ws.Range("H3").Value = 68.10198299999999
$ws.Range("I3").Value = 0.08615268874617349
$ws.Range("J3").Value = 0.08615268874617349
$ws.Range("O3").Value = 0.1537223653287423
$ws.Range("P3").Value = 0.1537223653287423
$ws.Range("Q3").Value = 250.939500713515
$ws.Range("R3").Value = 2258.455506421635
$ws.Range("S3").Value = 0.01324359509349271
$ws.Range("T3").Value = 0.01324359509349271
$ws.Range("G4").Value = 22.700661
$ws.Range("H4").Value = 68.10198299999999
$ws.Range("I4").Value = 0.08615268874617349
$ws.Range("J4").Value = 0.08615268874617349
$ws.Range("M4").Value = 30.561198
$ws.Range("N4").Value = 91.683594
$ws.Range("O4").Value = 0.4249882340167162
$ws.Range("P4").Value = 0.4249882340167161
$ws.Range("Q4").Value = 693.7593955518779
$ws.Range("R4").Value = 6243.834559966901
$ws.Range("S4").Value = 0.03661387904602809
$ws.Range("T4").Value = 0.03661387904602809
$ws.Range("G5").Value = 22.700661
$ws.Range("H5").Value = 68.10198299999999
$ws.Range("I5").Value = 0.08615268874617349
$ws.Range("J5").Value = 0.08615268874617349
$ws.Range("M5").Value = 29.75783666666667
$ws.Range("N5").Value = 89.27351
$ws.Range("O5").Value = 0.4138165805255589
$ws.Range("P5").Value = 0.4138165805255589
$ws.Range("Q5").Value = 675.5225622633699
$ws.Range("R5").Value = 6079.703060370329
$ws.Range("S5").Value = 0.03565141106002431
$ws.Range("T5").Value = 0.03565141106002431
$ws.Range("I6").Value = 0.5030288587986086
$ws.Range("J6").Value = 0.5030288587986087
$ws.Range("M6").Value = 0.5373756666666667
$ws.Range("N6").Value = 1.612127
$ws.Range("O6").Value = 0.007472820128982582
$ws.Range("P6").Value = 0.007472820128982581
$ws.Range("Q6").Value = 71.22632966379945
$ws.Range("R6").Value = 641.036966974195
$ws.Range("S6").Value = 0.00375904418148938
$ws.Range("T6").Value = 0.00375904418148938
$ws.Range("I7").Value = 0.5030288587986086
$ws.Range("J7").Value = 0.5030288587986087
$ws.Range("O7").Value = 0.1537223653287423
$ws.Range("P7").Value = 0.1537223653287423
$ws.Range("S7").Value = 0.07732678600314005
$ws.Range("T7").Value = 0.07732678600314005
$ws.Range("I8").Value = 0.5030288587986086
$ws.Range("J8").Value = 0.5030288587986087
$ws.Range("M8").Value = 30.561198
$ws.Range("N8").Value = 91.683594
$ws.Range("O8").Value = 0.4249882340167162
$ws.Range("P8").Value = 0.4249882340167161
$ws.Range("Q8").Value = 4050.72670515781
$ws.Range("R8").Value = 36456.54034642029
$ws.Range("S8").Value = 0.2137813463602648
$ws.Range("T8").Value = 0.2137813463602648
$ws.Range("I9").Value = 0.5030288587986086
$ws.Range("J9").Value = 0.5030288587986087
$ws.Range("M9").Value = 29.75783666666667
$ws.Range("N9").Value = 89.27351
$ws.Range("O9").Value = 0.4138165805255589
$ws.Range("P9").Value = 0.4138165805255589
$ws.Range("Q9").Value = 3944.245368698927
$ws.Range("R9").Value = 35498.20831829035
$ws.Range("S9").Value = 0.2081616822537144
$ws.Range("T9").Value = 0.2081616822537145
$ws.Range("G10").Value = 41.94534433333333
$ws.Range("H10").Value = 125.836033
$ws.Range("I10").Value = 0.159189381961201
$ws.Range("J10").Value = 0.159189381961201
$ws.Range("M10").Value = 0.5373756666666667
$ws.Range("N10").Value = 1.612127
$ws.Range("O10").Value = 0.007472820128982582
$ws.Range("P10").Value = 0.007472820128982581
$ws.Range("Q10").Value = 22.54040737468789
$ws.Range("R10").Value = 202.863666372191
$ws.Range("S10").Value = 0.00118959361783996
$ws.Range("T10").Value = 0.001189593617839959
$ws.Range("G11").Value = 41.94534433333333
$ws.Range("H11").Value = 125.836033
$ws.Range("I11").Value = 0.159189381961201
$ws.Range("J11").Value = 0.159189381961201
$ws.Range("O11").Value = 0.1537223653287423
$ws.Range("P11").Value = 0.1537223653287423
$ws.Range("Q11").Value = 463.6756508659873
$ws.Range("R11").Value = 4173.080857793885
$ws.Range("S11").Value = 0.02447096833029645
$ws.Range("T11").Value = 0.02447096833029644
$ws.Range("G12").Value = 41.94534433333333
$ws.Range("H12").Value = 125.836033
$ws.Range("I12").Value = 0.159189381961201
$ws.Range("J12").Value = 0.159189381961201
$ws.Range("M12").Value = 30.561198
$ws.Range("N12").Value = 91.683594
$ws.Range("O12").Value = 0.4249882340167162
$ws.Range("P12").Value = 0.4249882340167161
$ws.Range("Q12").Value = 1281.899973349178
$ws.Range("R12").Value = 11537.0997601426
$ws.Range("S12").Value = 0.06765361431390331
$ws.Range("T12").Value = 0.06765361431390331
$ws.Range("G13").Value = 41.94534433333333
$ws.Range("H13").Value = 125.836033
$ws.Range("I13").Value = 0.159189381961201
$ws.Range("J13").Value = 0.159189381961201
$ws.Range("M13").Value = 29.75783666666667
$ws.Range("N13").Value = 89.27351
$ws.Range("O13").Value = 0.4138165805255589
$ws.Range("P13").Value = 0.4138165805255589
$ws.Range("Q13").Value = 1248.202705598425
$ws.Range("R13").Value = 11233.82435038583
$ws.Range("S13").Value = 0.06587520569916129
$ws.Range("T13").Value = 0.06587520569916129
$ws.Range("G14").Value = 66.302588
$ws.Range("H14").Value = 198.907764
$ws.Range("I14").Value = 0.2516290704940168
$ws.Range("J14").Value = 0.2516290704940168
$ws.Range("M14").Value = 0.5373756666666667
$ws.Range("N14").Value = 1.612127
$ws.Range("O14").Value = 0.007472820128982582
$ws.Range("P14").Value = 0.007472820128982581
$ws.Range("Q14").Value = 35.62939742822534
$ws.Range("R14").Value = 320.664576854028
$ws.Range("S14").Value = 0.001880378783024866
$ws.Range("T14").Value = 0.001880378783024866
$ws.Range("G15").Value = 66.302588
$ws.Range("H15").Value = 198.907764
$ws.Range("I15").Value = 0.2516290704940168
$ws.Range("J15").Value = 0.2516290704940168
$ws.Range("O15").Value = 0.1537223653287423
$ws.Range("P15").Value = 0.1537223653287423
$ws.Range("Q15").Value = 732.9274829809534
$ws.Range("R15").Value = 6596.34734682858
$ws.Range("S15").Value = 0.03868101590181312
$ws.Range("T15").Value = 0.03868101590181311
$ws.Range("G16").Value = 66.302588
$ws.Range("H16").Value = 198.907764
$ws.Range("I16").Value = 0.2516290704940168
$ws.Range("J16").Value = 0.2516290704940168
$ws.Range("M16").Value = 30.561198
$ws.Range("N16").Value = 91.683594
$ws.Range("O16").Value = 0.4249882340167162
$ws.Range("P16").Value = 0.4249882340167161
$ws.Range("Q16").Value = 2026.286519780424
$ws.Range("R16").Value = 18236.57867802382
$ws.Range("S16").Value = 0.10693939429652
$ws.Range("T16").Value = 0.10693939429652
$ws.Range("G17").Value = 66.302588
$ws.Range("H17").Value = 198.907764
$ws.Range("I17").Value = 0.2516290704940168
$ws.Range("J17").Value = 0.2516290704940168
$ws.Range("M17").Value = 29.75783666666667
$ws.Range("N17").Value = 89.27351
$ws.Range("O17").Value = 0.4138165805255589
$ws.Range("P17").Value = 0.4138165805255589
$ws.Range("Q17").Value = 1973.021584281293
$ws.Range("R17").Value = 17757.19425853164
$ws.Range("S17").Value = 0.1041282815126589
$ws.Range("T17").Value = 0.1041282815126589
